# Code clean up for grade title:
# The "Grade Title" column (D) on the "Levels" sheet contained the same
# constant value "AAAS" in every data row (D3:D64). Remove that redundant
# value from the column (clearing the cells, but keeping their formatting),
# which also drops the now-unused "AAAS" entry from the shared string table
# on save. Finally, leave D3:D64 selected, matching the saved selection in
# the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Levels")

$ws.Range("D3:D64").ClearContents()
$ws.Range("D3:D64").Select()
